# ---------------------------------------------------------------------------
# analisis final del desafio #2
#
# 1) Remove the leading "Realizado por " paragraph.
# 2) "...Por ultimo, la funcion de despliegue..." -> "...Por ultimo, el
#    metodo de despliegue..." (also clears the spell-check proofErr marks
#    that wrapped the old word "funcion").
# 3) Append the new "CALCULO DE GASTO EN MEMORIA" paragraph content into the
#    empty paragraph just before the very last paragraph of the document.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1) Drop the "Realizado por " paragraph -------------------------------
$first = $d.Paragraphs(1)
if ($first.Range.Text.StartsWith("Realizado por")) {
    $first.Range.Delete()
}

# --- 2) "la funcion de despliegue" -> "el metodo de despliegue" -----------
$anchor = $d.Content.Duplicate
$anchor.Find.Execute("la funcion de despliegue", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$matchStart = $anchor.Start

# Delete "la funcion " (11 chars, through the trailing space) in one shot so
# the spell-check proofErr anchors tied to the removed run disappear too.
$killRange = $d.Range($matchStart, $matchStart + 11)
$killRange.Delete()

# Re-insert the replacement wording (with the trailing space put back).
$insPoint = $d.Range($matchStart, $matchStart)
$insPoint.InsertAfter("el método ")

# --- 3) New "CALCULO DE GASTO EN MEMORIA" paragraph content ---------------
$newText = "CÁLCULO DE GASTO EN MEMORIA: para poder calcular el valor aproximado del gasto de la memoria durante todo el programa, se utilizó el método de estimar un valor promedio superior y un promedio inferior. Ambos promedios (superior e inferior) se sacaron teniendo en cuenta cuanto costaba cada objeto por la cantidad de objetos existentes en cada clase correspondiente, al final todos estos valores, se sumaron, sacando como resultado el consumo general de todo el programa.  "

# Locate the empty paragraph that immediately precedes the document's final
# (also empty) paragraph.
$count = $d.Paragraphs.Count
$target = $d.Paragraphs($count - 1)

# Borrow run formatting (rFonts + lang) from a run elsewhere in the document
# that already carries exactly the formatting this paragraph's runs use, so
# the freshly typed text picks up matching <w:rPr> (rFonts + lang) instead
# of landing with no run properties at all.
$donor = $d.Content.Duplicate
$donor.Find.Execute("reto que se debe afrontar", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$insertStart = $target.Range.Start
$stamp = $d.Range($insertStart, $insertStart)
$stamp.FormattedText = $donor.FormattedText

$stamped = $d.Range($insertStart, $insertStart + $donor.Text.Length)
$stamped.Text = $newText
